$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34 - this shifts the existing rows 34-64 down to 35-65
$ws.Rows.Item(34).Insert()

# Populate the new row 34 (a duplicate of the former row 34, now at row 35,
# but with an updated date and volume)
$ws.Cells.Item(34, 1).Value = 4
$ws.Cells.Item(34, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value = "Los Lagos"
$ws.Cells.Item(34, 4).Value = [DateTime]"2022-11-30"
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
$ws.Cells.Item(34, 5).Value = 10
$ws.Cells.Item(34, 6).Value = 300000000
$ws.Cells.Item(34, 7).Value = "Espárragos"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 400
$ws.Cells.Item(34, 11).Value = 1800
$ws.Cells.Item(34, 12).Value = 1800
$ws.Cells.Item(34, 13).Value = 1800
$ws.Cells.Item(34, 14).Value = "$/kilo"
$ws.Cells.Item(34, 15).Value = "Provincia de Linares"
$ws.Cells.Item(34, 16).Value = 1800
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
